$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 24 updates
$ws.Range("E24").Value = "2026-02-12T16:15:57.750784+00:00"
$ws.Range("H24").Value = 7
$ws.Range("L24").Value = "[63, 31910, 19424, 30964, 30729, 19392, 29357]"

# Row 25 updates
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = "2026-02-12T16:16:11.612753+00:00"
$ws.Range("H25").Value = 2
$ws.Range("I25").Value = 1
$ws.Range("L25").Value = "[3, 2]"
$ws.Range("M25").Value = "[4]"
